$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# The second data row's Login_id ("mojizabidi98") is capitalized to match the
# naming convention used elsewhere (MOJIZABIDI98).
$ws.Range("E2").Value = "MOJIZABIDI98"

# The saved view had scrolled right (topLeftCell="Y1") with AA1:AB1 selected;
# reset the view back to the top-left and leave the selection on F8.
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("F8").Select()
